$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new topic "Big Data: What is it, how do I work with it?" into the
# schedule at D23, pushing the existing topics in D23:D27 down one row into
# D24:D28 (D28 previously had no topic cell).
#
# First, copy the formatting (style) of D23 down into D28 so the new cell
# picks up the same style used by the rest of this column (s="3"), matching
# the existing neighboring cells, then overwrite all six cell values.
$ws.Range("D23").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("D23").Value = "Big Data: What is it, how do I work with it?"
$ws.Range("D24").Value = "Plaintext Manipulation"
$ws.Range("D25").Value = "Data Science: Questions"
$ws.Range("D26").Value = "Data Science: Backwards Design"
$ws.Range("D27").Value = "Data Science: Tool Selection"
$ws.Range("D28").Value = "Project Proposal Workshopping"

# Reflect where the author's selection ended up after making this edit.
$ws.Range("D23").Select() | Out-Null
